$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1905
$ws1.Range("F3").Value = 1529
$ws1.Range("F5").Value = 791
$ws1.Range("F6").Value = 13408
$ws1.Range("F7").Value = 13273
$ws1.Range("F9").Value = 782
$ws1.Range("F13").Value = 693
$ws1.Range("F14").Value = 2103
$ws1.Range("F15").Value = 6
$ws1.Range("F21").Value = 282
$ws1.Range("F23").Value = 431
$ws1.Range("F24").Value = 767
$ws1.Range("F25").Value = 26

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 36
$ws2.Range("F5").Value = 134
$ws2.Range("F7").Value = 130
$ws2.Range("F8").Value = 135
$ws2.Range("F11").Value = 37

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 198
$ws3.Range("F3").Value = 54

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 198
$ws4.Range("F3").Value = 1905
$ws4.Range("F4").Value = 1529
$ws4.Range("F7").Value = 791
$ws4.Range("F8").Value = 13408
$ws4.Range("F9").Value = 13273
$ws4.Range("F11").Value = 782
$ws4.Range("F15").Value = 693
$ws4.Range("F16").Value = 36
$ws4.Range("F18").Value = 2103
$ws4.Range("F19").Value = 6
$ws4.Range("F23").Value = 134
$ws4.Range("F26").Value = 54
$ws4.Range("F28").Value = 282
$ws4.Range("F30").Value = 431
$ws4.Range("F31").Value = 767
$ws4.Range("F32").Value = 130
$ws4.Range("F33").Value = 136
$ws4.Range("F36").Value = 26
$ws4.Range("F37").Value = 37

$wb.Save()
